$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename (A1:D1) ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Title-case fix for Spanish connector words (de/del/la/las/los/el/y) in state/municipality names ---
# and one typo fix (MonteMorelos -> Montemorelos)
$textChanges = @(
    @{Cell='B8'; Value='Pabellón De Arteaga'},
    @{Cell='B9'; Value='Rincón De Romos'},
    @{Cell='B10'; Value='San Francisco De Los Romo'},
    @{Cell='B11'; Value='San José De Gracia'},
    @{Cell='B16'; Value='Playas De Rosarito'},
    @{Cell='B39'; Value='Amatenango De La Frontera'},
    @{Cell='B43'; Value='Benemérito De Las Américas'},
    @{Cell='B52'; Value='Chiapa De Corzo'},
    @{Cell='B57'; Value='Comitán De Domínguez'},
    @{Cell='B82'; Value='Mazapa De Madero'},
    @{Cell='B85'; Value='Montecristo De Guerrero'},
    @{Cell='B89'; Value='Ocozocoautla De Espinosa'},
    @{Cell='B97'; Value='Salto De Agua'},
    @{Cell='B98'; Value='San Cristóbal De Las Casas'},
    @{Cell='B141'; Value='Coyame Del Sotol'},
    @{Cell='B152'; Value='Guadalupe Y Calvo'},
    @{Cell='B155'; Value='Hidalgo Del Parral'},
    @{Cell='B178'; Value='San Francisco De Borja'},
    @{Cell='B179'; Value='San Francisco De Conchos'},
    @{Cell='B180'; Value='San Francisco Del Oro'},
    @{Cell='B188'; Value='Valle De Zaragoza'},
    @{Cell='B220'; Value='San Juan De Sabinas'},
    @{Cell='B237'; Value='Villa De Álvarez'},
    @{Cell='A239'; Value='Ciudad De México'},
    @{Cell='B243'; Value='Cuajimalpa De Morelos'},
    @{Cell='B258'; Value='Coneto De Comonfort'},
    @{Cell='B272'; Value='Nombre De Dios'},
    @{Cell='B276'; Value='Pánuco De Coronado'},
    @{Cell='B283'; Value='San Juan De Guadalupe'},
    @{Cell='B284'; Value='San Juan Del Río'},
    @{Cell='B285'; Value='San Luis Del Cordero'},
    @{Cell='B286'; Value='San Pedro Del Gallo'},
    @{Cell='A296'; Value='Estado De México'},
    @{Cell='B296'; Value='Acambay De Ruíz Castañeda'},
    @{Cell='B299'; Value='Almoloya De Alquisiras'},
    @{Cell='B300'; Value='Almoloya De Juárez'},
    @{Cell='B307'; Value='Atizapán De Zaragoza'},
    @{Cell='B315'; Value='Chapa De Mota'},
    @{Cell='B320'; Value='Coacalco De Berriozábal'},
    @{Cell='B327'; Value='Ecatepec De Morelos'},
    @{Cell='B335'; Value='Ixtapan De La Sal'},
    @{Cell='B336'; Value='Ixtapan Del Oro'},
    @{Cell='B353'; Value='Naucalpan De Juárez'},
    @{Cell='B367'; Value='San Antonio La Isla'},
    @{Cell='B368'; Value='San Felipe Del Progreso'},
    @{Cell='B370'; Value='San Simón De Guerrero'},
    @{Cell='B372'; Value='Soyaniquilpan De Juárez'},
    @{Cell='B382'; Value='Tenango Del Aire'},
    @{Cell='B383'; Value='Tenango Del Valle'},
    @{Cell='B397'; Value='Tlalnepantla De Baz'},
    @{Cell='B403'; Value='Valle De Bravo'},
    @{Cell='B404'; Value='Valle De Chalco Solidaridad'},
    @{Cell='B405'; Value='Villa De Allende'},
    @{Cell='B406'; Value='Villa Del Carbón'},
    @{Cell='B420'; Value='Apaseo El Alto'},
    @{Cell='B421'; Value='Apaseo El Grande'},
    @{Cell='B429'; Value='Dolores Hidalgo Cuna De La Independencia Nacional'},
    @{Cell='B433'; Value='Jaral Del Progreso'},
    @{Cell='B441'; Value='Purísima Del Rincón'},
    @{Cell='B445'; Value='San Diego De La Unión'},
    @{Cell='B447'; Value='San Francisco Del Rincón'},
    @{Cell='B449'; Value='San Luis De La Paz'},
    @{Cell='B451'; Value='Santa Cruz De Juventino Rosas'},
    @{Cell='B453'; Value='Silao De La Victoria'},
    @{Cell='B458'; Value='Valle De Santiago'},
    @{Cell='B464'; Value='Acapulco De Juárez'},
    @{Cell='B467'; Value='Ajuchitlán Del Progreso'},
    @{Cell='B468'; Value='Alcozauca De Guerrero'},
    @{Cell='B472'; Value='Atenango Del Río'},
    @{Cell='B473'; Value='Atlamajalcingo Del Monte'},
    @{Cell='B475'; Value='Atoyac De Álvarez'},
    @{Cell='B476'; Value='Ayutla De Los Libres'},
    @{Cell='B479'; Value='Buenavista De Cuéllar'},
    @{Cell='B480'; Value='Chilapa De Álvarez'},
    @{Cell='B481'; Value='Chilpancingo De Los Bravo'},
    @{Cell='B482'; Value='Coahuayutla De José María Izazaga'},
    @{Cell='B487'; Value='Coyuca De Benítez'},
    @{Cell='B488'; Value='Coyuca De Catalán'},
    @{Cell='B492'; Value='Cuetzala Del Progreso'},
    @{Cell='B493'; Value='Cutzamala De Pinzón'},
    @{Cell='B499'; Value='Huitzuco De Los Figueroa'},
    @{Cell='B500'; Value='Iguala De La Independencia'},
    @{Cell='B502'; Value='Ixcateopan De Cuauhtémoc'},
    @{Cell='B503'; Value='Zihuatanejo De Azueta'},
    @{Cell='B505'; Value='La Unión De Isidoro Montes De Oca'},
    @{Cell='B508'; Value='Mártir De Cuilapan'},
    @{Cell='B521'; Value='Taxco De Alarcón'},
    @{Cell='B523'; Value='Técpan De Galeana'},
    @{Cell='B525'; Value='Tepecoacuilco De Trujano'},
    @{Cell='B527'; Value='Tixtla De Guerrero'},
    @{Cell='B531'; Value='Tlalixtaquilla De Maldonado'},
    @{Cell='B532'; Value='Tlapa De Comonfort'},
    @{Cell='B544'; Value='Agua Blanca De Iturbide'},
    @{Cell='B551'; Value='Atotonilco De Tula'},
    @{Cell='B552'; Value='Atotonilco El Grande'},
    @{Cell='B558'; Value='Cuautepec De Hinojosa'},
    @{Cell='B564'; Value='Huasca De Ocampo'},
    @{Cell='B568'; Value='Huejutla De Reyes'},
    @{Cell='B571'; Value='Jacala De Ledezma'},
    @{Cell='B578'; Value='Mineral De La Reforma'},
    @{Cell='B579'; Value='Mineral Del Chico'},
    @{Cell='B580'; Value='Mineral Del Monte'},
    @{Cell='B581'; Value='Mixquiahuala De Juárez'},
    @{Cell='B582'; Value='Molango De Escamilla'},
    @{Cell='B584'; Value='Nopala De Villagrán'},
    @{Cell='B585'; Value='Omitlán De Juárez'},
    @{Cell='B586'; Value='Pachuca De Soto'},
    @{Cell='B589'; Value='Progreso De Obregón'},
    @{Cell='B595'; Value='Santiago De Anaya'},
    @{Cell='B596'; Value='Santiago Tulantepec De Lugo Guerrero'},
    @{Cell='B600'; Value='Tenango De Doria'},
    @{Cell='B602'; Value='Tepehuacán De Guerrero'},
    @{Cell='B603'; Value='Tepeji Del Río De Ocampo'},
    @{Cell='B606'; Value='Tezontepec De Aldama'},
    @{Cell='B615'; Value='Tula De Allende'},
    @{Cell='B616'; Value='Tulancingo De Bravo'},
    @{Cell='B617'; Value='Villa De Tezontepec'},
    @{Cell='B621'; Value='Zacualtipán De Ángeles'},
    @{Cell='B622'; Value='Zapotlán De Juárez'},
    @{Cell='B627'; Value='Acatlán De Juárez'},
    @{Cell='B628'; Value='Ahualulco De Mercado'},
    @{Cell='B633'; Value='Atemajac De Brizuela'},
    @{Cell='B636'; Value='Atotonilco El Alto'},
    @{Cell='B638'; Value='Autlán De Navarro'},
    @{Cell='B644'; Value='Cañadas De Obregón'},
    @{Cell='B650'; Value='Cuautitlán De García Barragán'},
    @{Cell='B658'; Value='Encarnación De Díaz'},
    @{Cell='B665'; Value='Huejuquilla El Alto'},
    @{Cell='B666'; Value='Ixtlahuacán De Los Membrillos'},
    @{Cell='B667'; Value='Ixtlahuacán Del Río'},
    @{Cell='B671'; Value='Jilotlán De Los Dolores'},
    @{Cell='B677'; Value='La Manzanilla De La Paz'},
    @{Cell='B678'; Value='Lagos De Moreno'},
    @{Cell='B686'; Value='Ojuelos De Jalisco'},
    @{Cell='B691'; Value='San Cristóbal De La Barranca'},
    @{Cell='B692'; Value='San Diego De Alejandría'},
    @{Cell='B694'; Value='San Juan De Los Lagos'},
    @{Cell='B695'; Value='San Juanito De Escobedo'},
    @{Cell='B698'; Value='San Martín De Bolaños'},
    @{Cell='B700'; Value='San Miguel El Alto'},
    @{Cell='B701'; Value='San Sebastián Del Oeste'},
    @{Cell='B702'; Value='Santa María De Los Ángeles'},
    @{Cell='B703'; Value='Santa María Del Oro'},
    @{Cell='B706'; Value='Talpa De Allende'},
    @{Cell='B707'; Value='Tamazula De Gordiano'},
    @{Cell='B710'; Value='Techaluta De Montenegro'},
    @{Cell='B714'; Value='Teocuitatlán De Corona'},
    @{Cell='B715'; Value='Tepatitlán De Morelos'},
    @{Cell='B718'; Value='Tizapán El Alto'},
    @{Cell='B719'; Value='Tlajomulco De Zúñiga'},
    @{Cell='B730'; Value='Unión De San Antonio'},
    @{Cell='B731'; Value='Unión De Tula'},
    @{Cell='B732'; Value='Valle De Guadalupe'},
    @{Cell='B733'; Value='Valle De Juárez'},
    @{Cell='B738'; Value='Yahualica De González Gallo'},
    @{Cell='B739'; Value='Zacoalco De Torres'},
    @{Cell='B742'; Value='Zapotitlán De Vadillo'},
    @{Cell='B743'; Value='Zapotlán Del Rey'},
    @{Cell='B744'; Value='Zapotlán El Grande'},
    @{Cell='B770'; Value='Coalcomán De Vázquez Pallares'},
    @{Cell='B772'; Value='Cojumatlán De Régules'},
    @{Cell='B839'; Value='Tiquicheo De Nicolás Romero'},
    @{Cell='B865'; Value='Coatlán Del Río'},
    @{Cell='B873'; Value='Jonacatepec De Leandro Valle'},
    @{Cell='B877'; Value='Puente De Ixtla'},
    @{Cell='B883'; Value='Tetela Del Volcán'},
    @{Cell='B884'; Value='Tlaltizapán De Zapata'},
    @{Cell='B892'; Value='Zacualpan De Amilpas'},
    @{Cell='B896'; Value='Amatlán De Cañas'},
    @{Cell='B897'; Value='Bahía De Banderas'},
    @{Cell='B901'; Value='Ixtlán Del Río'},
    @{Cell='B908'; Value='Santa María Del Oro'},
    @{Cell='B926'; Value='Ciénega De Flores'},
    @{Cell='B944'; Value='Lampazos De Naranjo'},
    @{Cell='B951'; Value='Mier Y Noriega'},
    @{Cell='B960'; Value='San Nicolás De Los Garza'},
    @{Cell='B968'; Value='Acatlán De Pérez Figueroa'},
    @{Cell='B972'; Value='Ayoquezco De Aldama'},
    @{Cell='B977'; Value='Chalcatongo De Hidalgo'},
    @{Cell='B978'; Value='Ciénega De Zimatlán'},
    @{Cell='B981'; Value='Coicoyán De Las Flores'},
    @{Cell='B984'; Value='Constancia Del Rosario'},
    @{Cell='B987'; Value='Cuilápam De Guerrero'},
    @{Cell='B989'; Value='Fresnillo De Trujano'},
    @{Cell='B991'; Value='Guelatao De Juárez'},
    @{Cell='B992'; Value='Guevea De Humboldt'},
    @{Cell='B993'; Value='Heroica Ciudad De Ejutla De Crespo'},
    @{Cell='B994'; Value='Heroica Ciudad De Huajuapan De León'},
    @{Cell='B995'; Value='Heroica Ciudad De Tlaxiaco'},
    @{Cell='B996'; Value='Huautla De Jiménez'},
    @{Cell='B998'; Value='Ixtlán De Juárez'},
    @{Cell='B999'; Value='Heroica Ciudad De Juchitán De Zaragoza'},
    @{Cell='B1011'; Value='Mariscala De Juárez'},
    @{Cell='B1012'; Value='Mártires De Tacubaya'},
    @{Cell='B1014'; Value='Mazatlán Villa De Flores'},
    @{Cell='B1015'; Value='Miahuatlán De Porfirio Díaz'},
    @{Cell='B1016'; Value='Mixistlán De La Reforma'},
    @{Cell='B1018'; Value='Nejapa De Madero'},
    @{Cell='B1020'; Value='Oaxaca De Juárez'},
    @{Cell='B1021'; Value='Ocotlán De Morelos'},
    @{Cell='B1022'; Value='Pinotepa De Don Luis'},
    @{Cell='B1024'; Value='Putla Villa De Guerrero'},
    @{Cell='B1025'; Value='Reforma De Pineda'},
    @{Cell='B1027'; Value='Rojas De Cuauhtémoc'},
    @{Cell='B1041'; Value='San Antonino El Alto'},
    @{Cell='B1043'; Value='San Antonio De La Cal'},
    @{Cell='B1047'; Value='San Baltazar Yatzachi El Bajo'},
    @{Cell='B1055'; Value='San Dionisio Del Mar'},
    @{Cell='B1074'; Value='San José Del Peñasco'},
    @{Cell='B1075'; Value='San José Del Progreso'},
    @{Cell='B1081'; Value='San Juan Bautista Lo De Soto'},
    @{Cell='B1090'; Value='San Juan De Los Cués'},
    @{Cell='B1091'; Value='San Juan Del Estado'},
    @{Cell='B1116'; Value='San Martín De Los Cansecos'},
    @{Cell='B1127'; Value='San Miguel Del Puerto'},
    @{Cell='B1140'; Value='San Pablo Villa De Mitla'},
    @{Cell='B1144'; Value='San Pedro El Alto'},
    @{Cell='B1159'; Value='San Pedro Y San Pablo Ayutla'},
    @{Cell='B1169'; Value='Santa Ana Del Valle'},
    @{Cell='B1178'; Value='Santa Cruz Tacache De Mina'},
    @{Cell='B1182'; Value='Santa Inés De Zaragoza'},
    @{Cell='B1183'; Value='Santa Lucía Del Camino'},
    @{Cell='B1197'; Value='Santa María Jalapa Del Marqués'},
    @{Cell='B1252'; Value='Santo Domingo De Morelos'},
    @{Cell='B1272'; Value='Tamazulápam Del Espíritu Santo'},
    @{Cell='B1273'; Value='Tataltepec De Valdés'},
    @{Cell='B1274'; Value='Teococuilco De Marcos Pérez'},
    @{Cell='B1275'; Value='Teotitlán De Flores Magón'},
    @{Cell='B1276'; Value='Teotitlán Del Valle'},
    @{Cell='B1278'; Value='Tezoatlán De Segura Y Luna'},
    @{Cell='B1279'; Value='Tlacolula De Matamoros'},
    @{Cell='B1281'; Value='Totontepec Villa De Morelos'},
    @{Cell='B1282'; Value='Villa De Chilapa De Díaz'},
    @{Cell='B1283'; Value='Villa De Etla'},
    @{Cell='B1284'; Value='Villa De Tututepec De Melchor Ocampo'},
    @{Cell='B1285'; Value='Villa De Zaachila'},
    @{Cell='B1286'; Value='Villa Sola De Vega'},
    @{Cell='B1288'; Value='Zimatlán De Álvarez'},
    @{Cell='B1313'; Value='Ayotoxco De Guerrero'},
    @{Cell='B1318'; Value='Chalchicomula De Sesma'},
    @{Cell='B1339'; Value='Cuapiaxtla De Madero'},
    @{Cell='B1342'; Value='Cuayuca De Andrade'},
    @{Cell='B1343'; Value='Cuetzalan Del Progreso'},
    @{Cell='B1359'; Value='Huehuetlán El Chico'},
    @{Cell='B1360'; Value='Huehuetlán El Grande'},
    @{Cell='B1365'; Value='Huitzilan De Serdán'},
    @{Cell='B1366'; Value='Ixcamilpa De Guerrero'},
    @{Cell='B1369'; Value='Izúcar De Matamoros'},
    @{Cell='B1379'; Value='Los Reyes De Juárez'},
    @{Cell='B1380'; Value='Mazapiltepec De Juárez'},
    @{Cell='B1390'; Value='Palmar De Bravo'},
    @{Cell='B1399'; Value='San Diego La Mesa Tochimiltzingo'},
    @{Cell='B1410'; Value='San Nicolás De Los Ranchos'},
    @{Cell='B1414'; Value='San Salvador El Seco'},
    @{Cell='B1415'; Value='San Salvador El Verde'},
    @{Cell='B1423'; Value='Tecali De Herrera'},
    @{Cell='B1431'; Value='Tepanco De López'},
    @{Cell='B1432'; Value='Tepango De Rodríguez'},
    @{Cell='B1433'; Value='Tepatlaxco De Hidalgo'},
    @{Cell='B1438'; Value='Tepexi De Rodríguez'},
    @{Cell='B1440'; Value='Tetela De Ocampo'},
    @{Cell='B1441'; Value='Teteles De Avila Castillo'},
    @{Cell='B1446'; Value='Tlacotepec De Benito Juárez'},
    @{Cell='B1458'; Value='Tuzamapan De Galeana'},
    @{Cell='B1472'; Value='Zapotitlán De Méndez'},
    @{Cell='B1479'; Value='Amealco De Bonfil'},
    @{Cell='B1481'; Value='Cadereyta De Montes'},
    @{Cell='B1487'; Value='Jalpan De Serra'},
    @{Cell='B1488'; Value='Landa De Matamoros'},
    @{Cell='B1491'; Value='Pinal De Amoles'},
    @{Cell='B1494'; Value='San Juan Del Río'},
    @{Cell='B1508'; Value='Armadillo De Los Infante'},
    @{Cell='B1509'; Value='Axtla De Terrazas'},
    @{Cell='B1514'; Value='Cerro De San Pedro'},
    @{Cell='B1516'; Value='Ciudad Del Maíz'},
    @{Cell='B1527'; Value='Mexquitic De Carmona'},
    @{Cell='B1533'; Value='San Ciro De Acosta'},
    @{Cell='B1539'; Value='Santa María Del Río'},
    @{Cell='B1541'; Value='Soledad De Graciano Sánchez'},
    @{Cell='B1549'; Value='Tanquián De Escobedo'},
    @{Cell='B1553'; Value='Villa De Arista'},
    @{Cell='B1554'; Value='Villa De Arriaga'},
    @{Cell='B1555'; Value='Villa De Guadalupe'},
    @{Cell='B1556'; Value='Villa De La Paz'},
    @{Cell='B1557'; Value='Villa De Ramos'},
    @{Cell='B1558'; Value='Villa De Reyes'},
    @{Cell='B1600'; Value='Nacozari De García'},
    @{Cell='B1621'; Value='Jalpa De Méndez'},
    @{Cell='B1666'; Value='Soto La Marina'},
    @{Cell='B1674'; Value='Acuamanala De Miguel Hidalgo'},
    @{Cell='B1676'; Value='Amaxac De Guerrero'},
    @{Cell='B1677'; Value='Apetatitlán De Antonio Carvajal'},
    @{Cell='B1683'; Value='Contla De Juan Cuamatzi'},
    @{Cell='B1689'; Value='Ixtacuixtla De Mariano Matamoros'},
    @{Cell='B1692'; Value='Mazatecochco De José María Morelos'},
    @{Cell='B1693'; Value='Nanacamilpa De Mariano Arista'},
    @{Cell='B1697'; Value='San Pablo Del Monte'},
    @{Cell='B1701'; Value='Tepetitla De Lardizábal'},
    @{Cell='B1704'; Value='Tetla De La Solidaridad'},
    @{Cell='B1715'; Value='Ziltlaltépec De Trinidad Sánchez Santos'},
    @{Cell='B1725'; Value='Alto Lucero De Gutiérrez Barrios'},
    @{Cell='B1729'; Value='Amatlán De Los Reyes'},
    @{Cell='B1740'; Value='Boca Del Río'},
    @{Cell='B1742'; Value='Camarón De Tejeda'},
    @{Cell='B1746'; Value='Castillo De Teayo'},
    @{Cell='B1748'; Value='Cazones De Herrera'},
    @{Cell='B1756'; Value='Chinampa De Gorostiza'},
    @{Cell='B1769'; Value='Cosamaloapan De Carpio'},
    @{Cell='B1770'; Value='Cosautlán De Carvajal'},
    @{Cell='B1787'; Value='Hueyapan De Ocampo'},
    @{Cell='B1788'; Value='Huiloapan De Cuauhtémoc'},
    @{Cell='B1789'; Value='Ignacio De La Llave'},
    @{Cell='B1793'; Value='Ixhuatlán De Madero'},
    @{Cell='B1794'; Value='Ixhuatlán Del Café'},
    @{Cell='B1795'; Value='Ixhuatlán Del Sureste'},
    @{Cell='B1806'; Value='Juchique De Ferrer'},
    @{Cell='B1809'; Value='Landero Y Coss'},
    @{Cell='B1811'; Value='Las Vigas De Ramírez'},
    @{Cell='B1812'; Value='Lerdo De Tejada'},
    @{Cell='B1817'; Value='Martínez De La Torre'},
    @{Cell='B1820'; Value='Medellín De Bravo'},
    @{Cell='B1825'; Value='Nanchital De Lázaro Cárdenas Del Río'},
    @{Cell='B1836'; Value='Ozuluama De Mascareñas'},
    @{Cell='B1840'; Value='Paso De Ovejas'},
    @{Cell='B1841'; Value='Paso Del Macho'},
    @{Cell='B1845'; Value='Poza Rica De Hidalgo'},
    @{Cell='B1855'; Value='Sayula De Alemán'},
    @{Cell='B1859'; Value='Soledad De Doblado'},
    @{Cell='B1867'; Value='Tatahuicapan De Juárez'},
    @{Cell='B1900'; Value='Vega De Alatorre'},
    @{Cell='B1911'; Value='Zontecomatlán De López Y Fuentes'},
    @{Cell='B1912'; Value='Zozocolco De Hidalgo'},
    @{Cell='B1953'; Value='Cañitas De Felipe Pescador'},
    @{Cell='B1955'; Value='Concepción Del Oro'},
    @{Cell='B1957'; Value='El Plateado De Joaquín Amaro'},
    @{Cell='B1968'; Value='Jiménez Del Teul'},
    @{Cell='B1975'; Value='Mezquital Del Oro'},
    @{Cell='B1980'; Value='Moyahua De Estrada'},
    @{Cell='B1981'; Value='Nochistlán De Mejía'},
    @{Cell='B1982'; Value='Noria De Ángeles'},
    @{Cell='B1993'; Value='Teúl De González Ortega'},
    @{Cell='B1994'; Value='Tlaltenango De Sánchez Román'},
    @{Cell='B1996'; Value='Trinidad García De La Cadena'},
    @{Cell='B1999'; Value='Villa De Cos'},
    @{Cell='B953'; Value='Montemorelos'},
)

foreach ($chg in $textChanges) {
    $ws.Range($chg.Cell).Value = $chg.Value
}

# --- Tiny floating point recomputation (last-bit ULP differences) in pct_matriculas column ---
$floatChanges = @(
    @{Cell='D9'; Value=0.0009660038210817812},
    @{Cell='D18'; Value=0.0009516926533620512},
    @{Cell='D296'; Value=0.0009660038210817812},
    @{Cell='D446'; Value=0.0091305250051878},
    @{Cell='D447'; Value=0.000958848237221916},
    @{Cell='D497'; Value=0.0009874705726613765},
    @{Cell='D1769'; Value=0.0009516926533620512},
    @{Cell='D1903'; Value=0.0009230703179225908},
)

foreach ($chg in $floatChanges) {
    $ws.Range($chg.Cell).Value = $chg.Value
}

# --- Remove trailing footer/metadata rows (2008:2012), which also shrinks the sheet dimension to A1:D2006 ---
$ws.Range('A2008:D2012').EntireRow.Delete() | Out-Null

Write-Host "Edit complete"
